# Weekly cryptos list refresh (GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the literal text into the cell without Excel coercing
    # number-looking strings (e.g. "1.00", "66.596.00") into numerics,
    # then drop the temporary Text number-format so no explicit style
    # sticks to the cell (matches the source data which is unstyled).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "66.596.00"
$ws.Range("E2").Value = "  +4.54%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.444.25"
$ws.Range("E3").Value = "  +5.13%  "

# Row 4
$ws.Range("E4").Value = "  -0.30%  "

# Row 5
Set-TextValue $ws.Range("D5") "185.03"
$ws.Range("E5").Value = "  +4.60%  "

# Row 6
Set-TextValue $ws.Range("D6") "548.05"
$ws.Range("E6").Value = "  +4.97%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.613"
$ws.Range("E7").Value = "  +2.14%  "

# Row 8
Set-TextValue $ws.Range("D8") "3.441.81"
$ws.Range("E8").Value = "  +5.37%  "

# Row 9
$ws.Range("E9").Value = "  -0.18%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.638"
$ws.Range("E10").Value = "  +5.44%  "

# Row 11
Set-TextValue $ws.Range("D11") "56.34"
$ws.Range("E11").Value = "  -1.90%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.148"
$ws.Range("E12").Value = "  +12.43%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.0000275"
$ws.Range("E13").Value = "  +7.20%  "

# Row 14
Set-TextValue $ws.Range("D14") "9.41"
$ws.Range("E14").Value = "  +4.16%  "

# Row 15
Set-TextValue $ws.Range("D15") "3.957.76"
$ws.Range("E15").Value = "  +4.10%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D16") "3.418.12"
$ws.Range("E16").Value = "  +4.07%  "

# Row 17
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D17") "0.121"
$ws.Range("E17").Value = "  +3.53%  "

# Row 18
Set-TextValue $ws.Range("D18") "18.28"
$ws.Range("E18").Value = "  +5.49%  "

# Row 19
Set-TextValue $ws.Range("D19") "66.672.99"
$ws.Range("E19").Value = "  +4.70%  "

# Row 20
Set-TextValue $ws.Range("D20") "11.73"
$ws.Range("E20").Value = "  +6.34%  "

# Row 21
Set-TextValue $ws.Range("D21") "1.00"
$ws.Range("E21").Value = "  +5.43%  "

# Row 22
Set-TextValue $ws.Range("D22") "405.55"
$ws.Range("E22").Value = "  +9.12%  "

# Row 23
Set-TextValue $ws.Range("D23") "12.03"
$ws.Range("E23").Value = "  +9.34%  "

# Row 24
Set-TextValue $ws.Range("D24") "4.25"
$ws.Range("E24").Value = "  +9.86%  "

# Row 25
Set-TextValue $ws.Range("D25") "3.89"
$ws.Range("E25").Value = "  +3.66%  "

# Row 26
Set-TextValue $ws.Range("D26") "84.30"
$ws.Range("E26").Value = "  +5.27%  "

# Row 27
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D27") "2.90"
$ws.Range("E27").Value = "  +9.08%  "

# Row 28
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D28") "6.23"
$ws.Range("E28").Value = "  +1.35%  "

# Row 29
Set-TextValue $ws.Range("D29") "11.73"
$ws.Range("E29").Value = "  +3.77%  "

# Row 30
Set-TextValue $ws.Range("D30") "8.63"
$ws.Range("E30").Value = "  +4.05%  "

# Row 31
Set-TextValue $ws.Range("D31") "30.19"
$ws.Range("E31").Value = "  +5.30%  "

# Row 32
Set-TextValue $ws.Range("D32") "674.67"
$ws.Range("E32").Value = "  +5.47%  "

# Row 33
Set-TextValue $ws.Range("D33") "6.84"
$ws.Range("E33").Value = "  +3.98%  "

# Row 34
Set-TextValue $ws.Range("D34") "11.60"
$ws.Range("E34").Value = "  +3.74%  "

# Row 35
$ws.Range("E35").Value = "  +4.89%  "

# Row 36
Set-TextValue $ws.Range("D36") "59.10"
$ws.Range("E36").Value = "  +0.28%  "

# Row 37
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D37") "38.75"
$ws.Range("E37").Value = "  +6.61%  "

# Row 38
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D38") "0.0₃0822"
$ws.Range("E38").Value = "  +18.70%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.403"
$ws.Range("E39").Value = "  +4.07%  "

# Row 40
$ws.Range("E40").Value = "  +0.14%  "

# Row 41
Set-TextValue $ws.Range("D41") "2.81"
$ws.Range("E41").Value = "  +15.26%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D42") "3.35"
$ws.Range("E42").Value = "  +20.71%  "

# Row 43
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D43") "0.133"
$ws.Range("E43").Value = "  +7.79%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.999"
$ws.Range("E44").Value = "  -0.14%  "

# Row 45
Set-TextValue $ws.Range("D45") "3.037.43"
$ws.Range("E45").Value = "  +3.48%  "

# Row 46
Set-TextValue $ws.Range("D46") "2.95"
$ws.Range("E46").Value = "  +10.07%  "

# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D47") "0.0419"
$ws.Range("E47").Value = "  +5.95%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D48") "3.27"
$ws.Range("E48").Value = "  +7.76%  "

# Row 49
$ws.Range("E49").Value = "  +4.06%  "

# Row 50
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D50") "2.66"
$ws.Range("E50").Value = "  +8.95%  "

# Row 51
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D51") "8.74"
$ws.Range("E51").Value = "  +11.65%  "
